# mrbeast-analysis: rename "Statics" -> "Statistics" and make it the
# active sheet/tab, with E21 selected (matching the shipped workbook).

$wb = $excel.ActiveWorkbook

# Rename the "Statics" worksheet to "Statistics"
$ws = $wb.Worksheets.Item("Statics")
$ws.Name = "Statistics"

# Switch focus to the (now renamed) Statistics sheet and select E21,
# which becomes the active/visible tab in the saved workbook.
$ws.Activate()
$ws.Range("E21").Select()
